$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Disaster" column (I) with header + values for each data row
$ws.Range("I1").Value = "Disaster"
$ws.Range("I2").Value = 0.5
$ws.Range("I3").Value = 0.5
$ws.Range("I4").Value = 0.5

# I5 previously only carried inherited formatting (style s="4") with no
# content; clear that formatting before writing the value so the cell
# reverts to the default style, matching the target state.
$ws.Range("I5").ClearFormats()
$ws.Range("I5").Value = 0.5

# Update the active selection shown when the sheet is next opened.
$ws.Range("G10").Select() | Out-Null
